$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are numeric-looking strings that must stay as text.
# Temporarily force text format, assign the value, then restore the default style
# so the cell keeps no explicit style (matches original formatting).
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "29.031.55"
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.834.07"
$cell.Style = "Normal"
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.9983"
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "242.37"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.6270"
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.9996"
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.07586"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.2925"
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "22.60"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07719"
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.836.56"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.953"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.6663"
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.00001018"
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "82.73"
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "6.041"
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "29.041.61"
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "226.74"
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "12.35"
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.9989"
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "7.180"
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.9994"
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "158.37"
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "8.498"
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.1375"
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "17.92"
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.491"
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "4.104"
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "4.020"
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.192"
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.05218"
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.847"
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.7369"
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.140"
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.703"
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.239.00"
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.757"
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.01786"
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "6.336"
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.8964"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.9996"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "101.54"
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.00000000125"
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.981.16"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "64.28"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.5108"
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.4039"
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "8.851"
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.644"
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.05748"
$cell.Style = "Normal"

# Coin name / link / volume columns are plain text already; assign directly.
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  -5.21%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +2.22%  "
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("E15").Value = "  +16.88%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("E37").Value = "  -4.73%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("E44").Value = "  +3.03%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E50").Value = "  -6.05%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("E51").Value = "  -1.88%  "
